# Split the reply-to-fifth-comment text so the parenthetical note reflects
# that replies can now be included via the "include_replies" config option.
#
# Original comment text (single run, after the annotationRef run):
#   "This is a reply to the fifth comment. (Not included)"
#
# New text (three runs, after the annotationRef run):
#   "This is a reply to the fifth comment. ("
#   "Can be included by using \u201cinclude_replies\u201d option in the config"
#   ")"

$d = $word.ActiveDocument

$part1 = "This is a reply to the fifth comment. ("
$part2 = [string][char]0x201C + "include_replies" + [string][char]0x201D
$part3 = "Can be included by using " + $part2 + " option in the config"
$part4 = ")"

# Locate the reply comment ("This is a reply to the fifth comment...") among
# the document's comments rather than assuming a fixed index.
$target = $null
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $candidate = $d.Comments.Item($i)
    if ($candidate.Range.Text -like "This is a reply to the fifth comment*") {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    $target.Range.Text = $part1 + $part3 + $part4
}
